# Update column F (dSF) values on the active sheet to reflect the
# repulled/recalculated data from the commit "repull data, push all data,
# mean calculation".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = -4
    3  = -1
    4  = -1
    6  = -4
    7  = -2
    8  = 4
    9  = 1
    10 = 3
    11 = 7
    12 = -2
    14 = -2
    15 = 4
    16 = -8
    17 = 7
    18 = -2
    20 = 8
    21 = 0
    22 = -5
    23 = 3
    24 = 3
}

foreach ($row in $newValues.Keys) {
    $ws.Range("F$row").Value = $newValues[$row]
}
